$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("For the duration of the project, storage and backup of data will be ensured by the project manager.", `
               $true, $false, $false, $false, $false, $true, 1, $false, `
               "[storageintro]", 2)
